$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 100 (shifts old row 100 -> row 101, preserving its content)
$ws.Rows.Item(100).Insert()

$data = @"
1|Richard|Louis|Exclu
2|Fontaine|Sacha|Exclu
3|Garnier|Louise|Exclu
4|Roger|Lina|Red
5|Garnier|Noah|Exclu
6|David|Camille|Red
7|Garcia|Léo|Red
8|Bertrand|Paul|Red
9|Richard|Manon|Exclu
10|Garcia|Lucas|Passe
11|Chevalier|Jade|Red
12|Dubois|Emma|Red
13|Martin|Lucas|Passe
14|Vincent|Louis|Red
15|Martin|Paul|Exclu
16|Michel|Alice|Passe
17|Roux|Noah|Exclu
18|Robert|Alice|Exclu
19|Durand|Hugo|Passe
20|Fontaine|Jade|Passe
21|Fontaine|Nathan|Passe
22|Moulin|Noah|Passe
23|Thomas|Paul|Exclu
24|Durand|Camille|Red
25|Morin|Arthur|Exclu
26|Morin|Arthur|Red
27|Moulin|Camille|Passe
28|Robert|Raphaël|Red
29|Simon|Adam|Exclu
30|Robert|Louis|Red
31|Fournier|Louis|Passe
32|Martin|Léo|Exclu
33|David|Paul|Exclu
34|Moulin|Inès|Passe
35|Thomas|Emma|Exclu
36|Thomas|Noah|Red
37|Morin|Manon|Passe
38|Garcia|Arthur|Exclu
39|Morin|Louise|Passe
40|Durand|Raphaël|Passe
41|David|Zoé|Red
42|Fontaine|Arthur|Exclu
43|Laurent|Gabriel|Exclu
44|Vincent|Camille|Passe
45|Garnier|Jules|Exclu
46|Fontaine|Hugo|Passe
47|Vincent|Emma|Exclu
48|Moulin|Chloé|Red
49|Moulin|Louis|Passe
50|Fontaine|Léo|Red
51|Fournier|Manon|Passe
52|Simon|Noah|Passe
53|Roux|Sarah|Exclu
54|David|Léo|Passe
55|Roger|Louise|Red
56|Morin|Zoé|Passe
57|Richard|Noah|Passe
58|Morin|Alice|Exclu
59|Fournier|Paul|Passe
60|Roger|Sarah|Red
61|Vincent|Arthur|Passe
62|Fontaine|Léa|Exclu
63|Petit|Nathan|Passe
64|Bertrand|Sacha|Exclu
65|Laurent|Adam|Red
66|Vincent|Lucas|Exclu
67|Richard|Jade|Red
68|Roger|Lucas|Red
69|Petit|Jules|Exclu
70|Michel|Adam|Red
71|Leroy|Alice|Red
72|Morin|Lucas|Exclu
73|Fournier|Jules|Exclu
74|Simon|Léo|Passe
75|Martin|Alice|Red
76|Bernard|Raphaël|Red
77|Leroy|Louis|Passe
78|David|Léo|Red
79|Moulin|Alice|Passe
80|Durand|Lucas|Passe
81|Michel|Paul|Red
82|Morin|Hugo|Exclu
83|Garcia|Léo|Exclu
84|Vincent|Jules|Passe
85|Durand|Adam|Red
86|Bertrand|Arthur|Passe
87|Moulin|Hugo|Red
88|Morin|Alice|Passe
89|Moreau|Camille|Exclu
90|Leroy|Adam|Exclu
91|Leroy|Emma|Red
92|Martin|Gabriel|Exclu
93|Bertrand|Inès|Passe
94|Petit|Gabriel|Red
95|Bertrand|Emma|Red
96|Richard|Raphaël|Exclu
97|Petit|Noah|Red
98|Thomas|Arthur|Passe
99|Fontaine|Hugo|Exclu
100|Vincent|Jules|Exclu
"@

$lines = $data -split "`n"

$rowIndex = 2
foreach ($line in $lines) {
    $line = $line.Trim()
    if ($line.Length -eq 0) { continue }
    $parts = $line -split '\|'
    $ws.Cells.Item($rowIndex, 1).Value = [int]$parts[0]
    $ws.Cells.Item($rowIndex, 2).Value = $parts[1]
    $ws.Cells.Item($rowIndex, 3).Value = $parts[2]
    $ws.Cells.Item($rowIndex, 18).Value = $parts[3]
    $rowIndex++
}

Write-Output "rows written: $($rowIndex - 2)"
